$d = $word.ActiveDocument

# Update the date/weekday heading line (unique text in the document).
$d.Content.Find.Execute("2026-01-10 Saturday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2026-01-11 Sunday", 2)

# Update the division drill table. The table has 20 physical rows (content
# rows interleaved with blank spacer rows) and 5 columns; only physical
# rows 1, 5, 9, 13, 17 hold the visible problems. Some values repeat
# ("63÷7=", "92÷6=") so we address cells directly by position instead of
# using a global text Find/Replace (which would be ambiguous / could
# clobber a just-written replacement).
$t = $d.Tables.Item(1)

$grid = @(
    @("10÷5=", "63÷7=", "98÷3=", "14÷4=", "58÷4="),
    @("25÷6=", "49÷6=", "13÷4=", "51÷6=", "74÷3="),
    @("60÷5=", "26÷9=", "56÷4=", "44÷3=", "26÷4="),
    @("30÷9=", "93÷9=", "69÷7=", "98÷3=", "36÷4="),
    @("29÷8=", "34÷8=", "78÷7=", "80÷8=", "32÷3=")
)

$physRows = @(1, 5, 9, 13, 17)

for ($i = 0; $i -lt $physRows.Length; $i++) {
    $r = $physRows[$i]
    $rowValues = $grid[$i]
    for ($c = 1; $c -le 5; $c++) {
        $t.Cell($r, $c).Range.Text = $rowValues[$c - 1]
    }
}
